$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''306.93'
$ws.Range("E2").Value = '''-3.52%'
$ws.Range("D3").Value = '''40.43'
$ws.Range("E3").Value = '''-3.74%'
$ws.Range("D4").Value = '''5.031'
$ws.Range("E4").Value = '''-3.13%'
$ws.Range("D5").Value = '''0.07604'
$ws.Range("E5").Value = '''-6.44%'
$ws.Range("D6").Value = '''4.243'
$ws.Range("E6").Value = '''-2.90%'
$ws.Range("D7").Value = '''1.593'
$ws.Range("E7").Value = '''-9.47%'
$ws.Range("D8").Value = '''0.9078'
$ws.Range("E8").Value = '''-2.55%'
$ws.Range("D9").Value = '''0.09991'
$ws.Range("E9").Value = '''-10.63%'
$ws.Range("D10").Value = '''0.1749'
$ws.Range("E10").Value = '''-6.00%'
$ws.Range("D11").Value = '''0.09014'
$ws.Range("E11").Value = '''-2.95%'
$ws.Range("D12").Value = '''0.04337'
$ws.Range("E12").Value = '''-5.32%'
$ws.Range("D13").Value = '''0.1055'
$ws.Range("E13").Value = '''-0.12%'
$ws.Range("D14").Value = '''0.001229'
$ws.Range("E14").Value = '''-3.73%'
$ws.Range("D15").Value = '''0.005855'
$ws.Range("E15").Value = '''1.75%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.370'
$ws.Range("E16").Value = '''0.62%'
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D17").Value = '''2.444'
$ws.Range("E17").Value = '''-3.44%'
$ws.Range("B18").Value = 'BitpandaEcosystemToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D18").Value = '''0.3279'
$ws.Range("E18").Value = '''-2.90%'
$ws.Range("B19").Value = 'MCDex'
$ws.Range("C19").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D19").Value = '''6.885'
$ws.Range("E19").Value = '''-7.12%'
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").Value = '''0.1351'
$ws.Range("E20").Value = '''-2.28%'
$ws.Range("B21").Value = 'ZBToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D21").Value = '''0.2850'
$ws.Range("E21").Value = '''9.61%'
$ws.Range("B22").Value = 'CoinExToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D22").Value = '''0.04171'
$ws.Range("E22").Value = '''-0.12%'
$ws.Range("B23").Value = 'BitKan'
$ws.Range("C23").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D23").Value = '''0.001218'
$ws.Range("E23").Value = '''-2.01%'
$ws.Range("B24").Value = 'HotbitToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D24").Value = '''0.004062'
$ws.Range("E24").Value = '''-4.56%'
$ws.Range("B25").Value = 'NitroEx'
$ws.Range("C25").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D25").Value = '''0.0001305'
$ws.Range("E25").Value = '''6.74%'
$ws.Range("B26").Value = 'Spectre.aiUtilityToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut'
$ws.Range("D26").Value = '--'
$ws.Range("E26").Value = '--%'
$ws.Range("B27").Value = 'LegolasExchange'
$ws.Range("C27").Value = 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo'
$ws.Range("B28").Value = 'BitZToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz'
$ws.Range("B29").Value = 'Birake'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir'
$ws.Range("B30").Value = 'NashExchange'
$ws.Range("C30").Value = 'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex'
$ws.Range("B31").Value = 'AAXToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
$ws.Range("B32").Value = 'CenX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx'
$ws.Range("B33").Value = 'BNIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix'
$ws.Range("B34").Value = 'UpBots'
$ws.Range("C34").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D38").Value = '''0.02420'
$ws.Range("E38").Value = '''-6.92%'
$ws.Range("E39").Value = '''-6.77%'
$ws.Range("D40").Value = '''0.007855'
$ws.Range("E40").Value = '''-2.50%'
$ws.Range("D41").Value = '''0.1304'
$ws.Range("E41").Value = '''-6.45%'
$ws.Range("D42").Value = '''0.007107'
$ws.Range("E42").Value = '''8.84%'
$ws.Range("D43").Value = '''0.001954'
$ws.Range("E43").Value = '''-6.19%'
$ws.Range("D44").Value = '''0.008367'
$ws.Range("E44").Value = '''1.84%'
$ws.Range("D45").Value = '''0.3318'
$ws.Range("E45").Value = '''-4.59%'
$ws.Range("D46").Value = '''0.00006450'
$ws.Range("E46").Value = '''-4.53%'
$ws.Range("E47").Value = '''0.01%'
$ws.Range("E48").Value = '''-26.89%'
$ws.Range("D49").Value = '''0.005336'
$ws.Range("E49").Value = '''57.38%'
$ws.Range("D50").Value = '''0.00002104'
$ws.Range("E50").Value = '''0.01%'
$ws.Range("D51").Value = '''0.0002004'
$ws.Range("E51").Value = '''0.01%'
